# Target sheet: "C330 comparisons" (3rd worksheet / sheetId=3 / r:id=rId3)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("C330 comparisons")
$ws.Activate()

# Column header labels used throughout this sheet's repeated blocks (A..Q)
$headerVals = @(
    "Year",
    " tot in HRUs reaches and reservoirs at end of last year (mm H2O)",
    " Precip (mm H2O)",
    " GW pumping (mm H2O)",
    " High Cascades groundwater contribution mm H2O",
    " from outside the basin (mm H2O)",
    " water added by FlowModel (mm)",
    " to outside the basin (mm H2O)",
    " AET (mm H2O)",
    " SNOW_EVAP (mm H2O)",
    " basin discharge (mm H2O)",
    " tot in HRUs reaches and reservoirs at end of this year (mm H2O)",
    " irrigation (ac-ft)",
    " municipal and rural domestic (ac-ft)",
    " mass balance discrepancy (mm H2O)",
    " mass balance discrepancy (fraction)",
    " weather year"
)
$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q")

# ---------------------------------------------------------------------
# 1) Row 43 (existing "C335+" block data row): add the per-column number
#    formatting that the rest of the sheet's data rows already carry
#    (columns B-L -> 0.00, M-N -> integer, O -> 0.00, P -> 0.000000).
# ---------------------------------------------------------------------
$ws.Range("B43:L43").NumberFormat = "0.00"
$ws.Range("M43:N43").NumberFormat = "0"
$ws.Range("O43").NumberFormat = "0.00"
$ws.Range("P43").NumberFormat = "0.000000"

# ---------------------------------------------------------------------
# 2) Row 45: blank structural separator row (cells B..P only, no value).
# ---------------------------------------------------------------------
$ws.Range("B45:P45").Style = "Normal"

# ---------------------------------------------------------------------
# 3) Row 46: new block label "C339"
# ---------------------------------------------------------------------
$ws.Range("A46").Value = "C339"

# ---------------------------------------------------------------------
# 4) Row 47: header row for the "C339" block
# ---------------------------------------------------------------------
for ($i = 0; $i -lt $headerVals.Length; $i++) {
    $ws.Range($cols[$i] + "47").Value = $headerVals[$i]
}
$ws.Range("A47:Q47").Style = "Normal"

# ---------------------------------------------------------------------
# 5) Row 48: data row for the "C339" block
# ---------------------------------------------------------------------
$ws.Range("A48").Value = 2010
$ws.Range("B48").Value = 1284.0238039999999
$ws.Range("C48").Value = 1990.4650879999999
$ws.Range("D48").Value = 1.4464170000000001
$ws.Range("E48").Value = 270.24752799999999
$ws.Range("F48").Value = 10.610913999999999
$ws.Range("G48").Value = 4.9719860000000002
$ws.Range("H48").Value = 8.8404570000000007
$ws.Range("I48").Value = 755.73742700000003
$ws.Range("J48").Value = 93.234084999999993
$ws.Range("K48").Value = 1371.6883539999999
$ws.Range("L48").Value = 1333.3901370000001
$ws.Range("M48").Value = 8273.0849610000005
$ws.Range("N48").Value = 29450.638672000001
$ws.Range("O48").Value = 1.124722
$ws.Range("P48").Value = 0.00031599999999999998
$ws.Range("Q48").Value = 2010

$ws.Range("B48:L48").NumberFormat = "0.00"
$ws.Range("M48:N48").NumberFormat = "0"
$ws.Range("O48").NumberFormat = "0.00"
$ws.Range("P48").NumberFormat = "0.000000"

# ---------------------------------------------------------------------
# 6) Row 50: new block label "C340" (row 49 intentionally left blank,
#    matching this sheet's existing spacer convention between blocks)
# ---------------------------------------------------------------------
$ws.Range("A50").Value = "C340"

# ---------------------------------------------------------------------
# 7) Row 51: header row for the "C340" block
# ---------------------------------------------------------------------
for ($i = 0; $i -lt $headerVals.Length; $i++) {
    $ws.Range($cols[$i] + "51").Value = $headerVals[$i]
}
$ws.Range("A51:Q51").Style = "Normal"

# ---------------------------------------------------------------------
# 8) Row 52: data row for the "C340" block
# ---------------------------------------------------------------------
$ws.Range("A52").Value = 2010
$ws.Range("B52").Value = 1284.0238039999999
$ws.Range("C52").Value = 1990.4650879999999
$ws.Range("D52").Value = 1.4464170000000001
$ws.Range("E52").Value = 270.24752799999999
$ws.Range("F52").Value = 10.610913999999999
$ws.Range("G52").Value = 4.9719850000000001
$ws.Range("H52").Value = 8.8404570000000007
$ws.Range("I52").Value = 755.73742700000003
$ws.Range("J52").Value = 93.234084999999993
$ws.Range("K52").Value = 1371.6883539999999
$ws.Range("L52").Value = 1333.3901370000001
$ws.Range("M52").Value = 8273.0849610000005
$ws.Range("N52").Value = 29450.638672000001
$ws.Range("O52").Value = 1.1247229999999999
$ws.Range("P52").Value = 0.00031599999999999998
$ws.Range("Q52").Value = 2010

$ws.Range("B52:L52").NumberFormat = "0.00"
$ws.Range("M52:N52").NumberFormat = "0"
$ws.Range("O52").NumberFormat = "0.00"
$ws.Range("P52").NumberFormat = "0.000000"

# ---------------------------------------------------------------------
# 9) Update the sheet view: selection moves to the new last header cell.
# ---------------------------------------------------------------------
$ws.Range("A51").Select()

Write-Output "C330 comparisons sheet updated through row 52"
